$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UI")
$ws.Activate()

# --- Update issue text for a few rows ---
$ws.Range("B2").Value = "When files are >1000(N), webserver is not responding"
$ws.Range("B14").Value = "When the directory/date folder entered is not available, there is no pop up/error saying ""no such directory exists""`n*Only list folder which are there in server"
$ws.Range("B17").Value = "Progress bar for data uploading"

# --- Apply the Status AutoFilter on C1:C18 (show In Progress / OnHold / blanks) ---
# Do this before filling in row 18 so the row (still blank at this point) is
# correctly folded into the "(Blanks)" bucket, matching the target filter def.
$rng = $ws.Range("C1:C18")
$criteria = @("In Progress", "OnHold", "")
$rng.AutoFilter(1, $criteria, 7)

# --- Fill in the previously-empty Status/Owner cells for the last row ---
$ws.Range("C18").Value = "In Progress"
$ws.Range("D18").Value = "Chakrapani"
foreach ($addr in @("C18","D18")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = 2
}

# --- Register the hidden _FilterDatabase name Excel creates for the filter ---
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=UI!`$C`$1:`$C`$18")
$fdb.Visible = $false

# --- Row 14 grew an extra wrapped line, so its row height increased ---
$ws.Rows.Item(14).RowHeight = 58

# --- Update the saved view: drop the old scroll/selection, select E9 ---
$ws.Range("E9").Select()
